$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quarterly database roll-forward: drop the oldest quarter
# ("فصل دوم منتهی به 1399/06"), shift every quarter column (E:N) one
# column to the left, and append the new quarter
# ("فصل چهارم منتهی به 1401/12") of data in column N.

# Row 8: Header row: quarter-end labels (E:N)
$values = @("فصل سوم منتهی به 1399/09", "فصل چهارم منتهی به 1399/12", "فصل اول منتهی به 1400/03", "فصل دوم منتهی به 1400/06", "فصل سوم منتهی به 1400/09", "فصل چهارم منتهی به 1400/12", "فصل اول منتهی به 1401/03", "فصل دوم منتهی به 1401/06", "فصل سوم منتهی به 1401/09", "فصل چهارم منتهی به 1401/12")
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(8, 5 + $i).Value = $values[$i]
}

# Row 10: هزینه حمل و نقل و انتقال
$values = @(2198444, 1649137, 1186321, 1181330, 807123, 833060, 1101716, 1725477, 1177496, 2423406)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(10, 5 + $i).Value = $values[$i]
}

# Row 12: حق العمل و کمیسیون فروش
$values = @(-2164, 151100, 2640, 2969, 537, 52913, 41116, 44305, -9481, -9120)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(12, 5 + $i).Value = $values[$i]
}

# Row 14: هزینه مواد مصرفی
$values = @(5016, 12994, 7729, 17780, 9110, 13381, 7343, 4486, 14136, 15465)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(14, 5 + $i).Value = $values[$i]
}

# Row 16: هزینه استهلاک
$values = @(1736, 2009, 1951, 1951, 9288, 6223, 3095, 3200, 3368, 3661)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(16, 5 + $i).Value = $values[$i]
}

# Row 17: هزینه حقوق و دستمزد
$values = @(59518, 113358, 107081, 138893, 110470, 117075, 157341, 217671, 165430, 247588)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(17, 5 + $i).Value = $values[$i]
}

# Row 19: سایر هزینه ها
$values = @(3252, 341660, 225408, 259128, 473908, 417254, 344107, 445591, 561951, 739935)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(19, 5 + $i).Value = $values[$i]
}

# Row 20: جمع
$values = @(2265802, 2270258, 1531130, 1602051, 1410436, 1439906, 1654718, 2440730, 1912900, 3420935)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(20, 5 + $i).Value = $values[$i]
}

# Row 24: Header row: quarter-end labels (E:N)
$values = @("فصل سوم منتهی به 1399/09", "فصل چهارم منتهی به 1399/12", "فصل اول منتهی به 1400/03", "فصل دوم منتهی به 1400/06", "فصل سوم منتهی به 1400/09", "فصل چهارم منتهی به 1400/12", "فصل اول منتهی به 1401/03", "فصل دوم منتهی به 1401/06", "فصل سوم منتهی به 1401/09", "فصل چهارم منتهی به 1401/12")
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(24, 5 + $i).Value = $values[$i]
}

# Row 26: تعداد پرسنل غیر تولیدی شرکت
$values = @(206, 217, 216, 217, 217, 220, 220, 220, 220, 220)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(26, 5 + $i).Value = $values[$i]
}

# Row 27: تعداد پرسنل تولیدی شرکت
$values = @(84, 84, 84, 84, 84, 85, 85, 85, 85, 85)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(27, 5 + $i).Value = $values[$i]
}
